# Añadiendo nombre de modelo
# Adds a "Modelo" column (F) describing the trained model pipeline,
# and refreshes a handful of metric values with their latest precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Model description text (identical for every data row) ---
$modelText = @"
MultiOutputRegressor(estimator=GridSearchCV(cv=5,
                                            estimator=Pipeline(steps=[('model',
                                                                       GradientBoostingRegressor())]),
                                            param_grid={'model__max_depth': [3,
                                                                             5,
                                                                             7],
                                                        'model__n_estimators': [50,
                                                                                100,
                                                                                150]},
                                            scoring='neg_mean_squared_error'))
"@

# --- Header cell F1: copy formatting from E1 (bordered/bold/centered header style), then set text ---
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Modelo"

# --- Fill F2:F14 with the model description text ---
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 6).Value = $modelText
    # Setting a multi-line value auto-expands the row height; restore the
    # original (default) row height/formatting so it matches the source.
    $ws.Rows.Item($row).AutoFit()
}

# --- Refresh a few metric values (MSE / MAE) with updated precision ---
$ws.Range("B3").Value = 0.1647362719327808
$ws.Range("D3").Value = 0.3215706312231011

$ws.Range("B4").Value = 0.2192987281846225
$ws.Range("D4").Value = 0.3777406617731509

$ws.Range("B5").Value = 0.3638702225807678
$ws.Range("D5").Value = 0.4682092966157618

$ws.Range("B7").Value = 0.07514644587374564
$ws.Range("D7").Value = 0.211919863475561

$ws.Range("B8").Value = 0.04215534119371416
$ws.Range("D8").Value = 0.136128825357167

$ws.Range("B9").Value = 0.07796894984218639
$ws.Range("D9").Value = 0.1911874935925048

$ws.Range("B11").Value = 0.0831172165082084
$ws.Range("D11").Value = 0.2072673588334908

$ws.Range("B13").Value = 0.0484517504435152
$ws.Range("D13").Value = 0.1526989685211048
